$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.250.11'
$ws.Range('E2').Value = '  +1.98%  '

$ws.Range('D3').Value = '3.317.50'
$ws.Range('E3').Value = '  +6.23%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.14'
$ws.Range('E5').Value = '  +1.13%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.55'
$ws.Range('E6').Value = '  +6.26%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').Value = '3.310.84'
$ws.Range('E8').Value = '  +6.07%  '

$ws.Range('E9').Value = '  +0.82%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.151'
$ws.Range('E10').Value = '  +3.01%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.48'
$ws.Range('E11').Value = '  +3.34%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.476'
$ws.Range('E12').Value = '  +3.61%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('E13').Value = '  +0.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.22'
$ws.Range('E14').Value = '  +3.61%  '

$ws.Range('D15').Value = '3.861.45'
$ws.Range('E15').Value = '  +6.07%  '

$ws.Range('D17').Value = '3.310.68'
$ws.Range('E17').Value = '  +5.86%  '

$ws.Range('D18').Value = '64.243.23'
$ws.Range('E18').Value = '  +1.97%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.94'
$ws.Range('E19').Value = '  +3.48%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '486.00'
$ws.Range('E20').Value = '  +2.20%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.43'
$ws.Range('E21').Value = '  +1.97%  '

$ws.Range('E22').Value = '  +7.39%  '

$ws.Range('E23').Value = '  +6.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.71'
$ws.Range('E24').Value = '  +5.13%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.22'
$ws.Range('E25').Value = '  -2.41%  '

$ws.Range('E26').Value = '  -0.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.79'
$ws.Range('E27').Value = '  +3.16%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.39'
$ws.Range('E28').Value = '  +4.93%  '

$ws.Range('E29').Value = '  +2.99%  '

$ws.Range('E30').Value = '  +0.01%  '

$ws.Range('E31').Value = '  +5.99%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.46'
$ws.Range('E32').Value = '  +4.40%  '

$ws.Range('E33').Value = '  +0.84%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.60'
$ws.Range('E34').Value = '  +2.48%  '

$ws.Range('E35').Value = '  +2.70%  '

$ws.Range('E36').Value = '  +3.53%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.24'
$ws.Range('E37').Value = '  +2.42%  '

$ws.Range('D38').Value = '0.0₃0750'
$ws.Range('E38').Value = '  +5.27%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0402'
$ws.Range('E39').Value = '  +2.98%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '432.79'
$ws.Range('E40').Value = '  +2.61%  '

$ws.Range('E41').Value = '  +4.61%  '

$ws.Range('D42').Value = '3.030.08'
$ws.Range('E42').Value = '  +5.96%  '

$ws.Range('E43').Value = '  +3.01%  '

$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.274'
$ws.Range('E44').Value = '  +6.51%  '

$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.110'
$ws.Range('E45').Value = '  -5.85%  '

$ws.Range('E46').Value = '  +8.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.54'
$ws.Range('E47').Value = '  +4.64%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.36'
$ws.Range('E48').Value = '  +3.28%  '

$ws.Range('E49').Value = '  +0.03%  '

$ws.Range('E50').Value = '  +1.68%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '123.93'
$ws.Range('E51').Value = '  +4.83%  '
